# Day4 of N8N learning
# Append a new journal row (row 5) documenting the "Mood playlisty fetcher"
# n8n workflow, matching the date-formatted entries already present.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data
$ws.Range("A5").Value = 45905
$ws.Range("A5").NumberFormat = "DD/MM/YY"
$ws.Range("B5").Value = "Mood playlisty fetcher"
$ws.Range("C5").Value = "fetches spotify playlist based on mood"
$ws.Range("D5").Value = "mood_chnager.json"

# Mirror the author's final selection position (moved a few rows further
# down after entering the new row).
$ws.Range("D10").Select() | Out-Null
